$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.00630148620626
$ws.Range("C2").Value = 9.142167332841243
$ws.Range("D2").Value = 14.84175391807416
$ws.Range("E2").Value = 16.24323577641535
$ws.Range("G2").Value = 29.27143134070009
$ws.Range("H2").Value = 14.26387854240084
$ws.Range("I2").Value = 19.39940973313773
$ws.Range("J2").Value = 9.257995831723399
$ws.Range("N2").Value = 16.51539280739487
$ws.Range("O2").Value = 21.81269242259754
$ws.Range("B3").Value = 13.44591234690282
$ws.Range("C3").Value = 8.644101103452142
$ws.Range("D3").Value = 14.77510407306477
$ws.Range("E3").Value = 16.1765548725051
$ws.Range("G3").Value = 29.20061464195817
$ws.Range("H3").Value = 14.30319562245235
$ws.Range("I3").Value = 19.4996097954302
$ws.Range("J3").Value = 9.264740972692078
$ws.Range("N3").Value = 16.55108749127059
$ws.Range("O3").Value = 21.85160147946895
$ws.Range("B4").Value = 13.0907732196211
$ws.Range("C4").Value = 8.323285427811443
$ws.Range("D4").Value = 14.7374645569726
$ws.Range("E4").Value = 16.13924414589047
$ws.Range("G4").Value = 29.16816003102913
$ws.Range("H4").Value = 14.33013901600449
$ws.Range("I4").Value = 19.5656577325697
$ws.Range("J4").Value = 9.270307462054827
$ws.Range("N4").Value = 16.57475533034225
$ws.Range("O4").Value = 21.88125219428225
$ws.Range("B5").Value = 12.94347732026318
$ws.Range("C5").Value = 8.188891969132058
$ws.Range("D5").Value = 14.72296384944497
$ws.Range("E5").Value = 16.12496438482972
$ws.Range("G5").Value = 29.15771385426775
$ws.Range("H5").Value = 14.34182216423614
$ws.Range("I5").Value = 19.59370833847844
$ws.Range("J5").Value = 9.272934381340752
$ws.Range("N5").Value = 16.58484118958395
$ws.Range("O5").Value = 21.89477935190897
$ws.Range("B6").Value = 12.91887002249256
$ws.Range("C6").Value = 8.166358846380181
$ws.Range("D6").Value = 14.72060694043558
$ws.Range("E6").Value = 16.12264940001853
$ws.Range("G6").Value = 29.15614725303095
$ws.Range("H6").Value = 14.3438045944764
$ws.Range("I6").Value = 19.59843463939814
$ws.Range("J6").Value = 9.273392238175406
$ws.Range("N6").Value = 16.58654259565411
$ws.Range("O6").Value = 21.89711262095302
$ws.Range("B7").Value = 13.08879686563856
$ws.Range("C7").Value = 8.321487596820553
$ws.Range("D7").Value = 14.73726558857482
$ws.Range("E7").Value = 16.13904780586392
$ws.Range("G7").Value = 29.16800789116779
$ws.Range("H7").Value = 14.33029373228238
$ws.Range("I7").Value = 19.56603143805194
$ws.Range("J7").Value = 9.270341437714935
$ws.Range("N7").Value = 16.57488956518288
$ws.Range("O7").Value = 21.88142878457665
$ws.Range("B8").Value = 13.8154909491428
$ws.Range("C8").Value = 8.973609024141966
$ws.Range("D8").Value = 14.81809903321325
$ws.Range("E8").Value = 16.21949769586953
$ws.Range("G8").Value = 29.24473012709266
$ws.Range("H8").Value = 14.27685268354031
$ws.Range("I8").Value = 19.43301826045362
$ws.Range("J8").Value = 9.260025962042304
$ws.Range("N8").Value = 16.52733728710481
$ws.Range("O8").Value = 21.82491029365841
$ws.Range("B9").Value = 15.14491352638535
$ws.Range("C9").Value = 10.12943330716353
$ws.Range("D9").Value = 15.00209057427476
$ws.Range("E9").Value = 16.40551708874919
$ws.Range("G9").Value = 29.48220982763124
$ws.Range("H9").Value = 14.19434478557109
$ws.Range("I9").Value = 19.20818230078014
$ws.Range("J9").Value = 9.251091661787845
$ws.Range("N9").Value = 16.44795245824524
$ws.Range("O9").Value = 21.75995204865617
$ws.Range("B10").Value = 16.05433656916108
$ws.Range("C10").Value = 10.89953114463716
$ws.Range("D10").Value = 15.15195620446312
$ws.Range("E10").Value = 16.55858986255313
$ws.Range("G10").Value = 29.70885132578004
$ws.Range("H10").Value = 14.14738401373901
$ws.Range("I10").Value = 19.06508353770382
$ws.Range("J10").Value = 9.251393057656971
$ws.Range("N10").Value = 16.39804143785009
$ws.Range("O10").Value = 21.74039065652627
$ws.Range("B11").Value = 16.45193541618388
$ws.Range("C11").Value = 11.2321003823563
$ws.Range("D11").Value = 15.22312286689083
$ws.Range("E11").Value = 16.6315915760071
$ws.Range("G11").Value = 29.82301383247187
$ws.Range("H11").Value = 14.12899942728973
$ws.Range("I11").Value = 19.00481178712047
$ws.Range("J11").Value = 9.253014815230017
$ws.Range("N11").Value = 16.37715429616105
$ws.Range("O11").Value = 21.73763705262748
$ws.Range("B12").Value = 16.60007757989359
$ws.Range("C12").Value = 11.35544673870198
$ws.Range("D12").Value = 15.250483879824
$ws.Range("E12").Value = 16.65970167257866
$ws.Range("G12").Value = 29.86780657339673
$ws.Range("H12").Value = 14.12246688210066
$ws.Range("I12").Value = 18.98268501002301
$ws.Range("J12").Value = 9.253841711215733
$ws.Range("N12").Value = 16.36950563269772
$ws.Range("O12").Value = 21.73747933234222
$ws.Range("B13").Value = 16.56828164269677
$ws.Range("C13").Value = 11.32899761367889
$ws.Range("D13").Value = 15.24457318639277
$ws.Range("E13").Value = 16.65362723977853
$ws.Range("G13").Value = 29.85809068637843
$ws.Range("H13").Value = 14.12385467362895
$ws.Range("I13").Value = 18.98741936649585
$ws.Range("J13").Value = 9.253654173852972
$ws.Range("N13").Value = 16.37114131644607
$ws.Range("O13").Value = 21.73747392445566
$ws.Range("B14").Value = 16.4641722039222
$ws.Range("C14").Value = 11.24230026270408
$ws.Range("D14").Value = 15.22536573062708
$ws.Range("E14").Value = 16.63389498098404
$ws.Range("G14").Value = 29.82666783342929
$ws.Range("H14").Value = 14.12845338137903
$ws.Range("I14").Value = 19.00297742068474
$ws.Range("J14").Value = 9.253078584689002
$ws.Range("N14").Value = 16.37651981131204
$ws.Range("O14").Value = 21.73760633449364
$ws.Range("B15").Value = 16.40008418464504
$ws.Range("C15").Value = 11.18885731191929
$ws.Range("D15").Value = 15.21365366826469
$ws.Range("E15").Value = 16.62186852028081
$ws.Range("G15").Value = 29.80762291576786
$ws.Range("H15").Value = 14.1313261596895
$ws.Range("I15").Value = 19.01259801725939
$ws.Range("J15").Value = 9.252753705732079
$ws.Range("N15").Value = 16.37984824874512
$ws.Range("O15").Value = 21.73780272315331
$ws.Range("B16").Value = 16.02801925616097
$ws.Range("C16").Value = 10.8774363284775
$ws.Range("D16").Value = 15.14736393838423
$ws.Range("E16").Value = 16.55388525094682
$ws.Range("G16").Value = 29.70161075036683
$ws.Range("H16").Value = 14.14864550863222
$ws.Range("I16").Value = 19.06911978135794
$ws.Range("J16").Value = 9.251316887822252
$ws.Range("N16").Value = 16.39944299157704
$ws.Range("O16").Value = 21.74069439586391
$ws.Range("B17").Value = 15.79556609282588
$ws.Range("C17").Value = 10.68181331273315
$ws.Range("D17").Value = 15.1074507454971
$ws.Range("E17").Value = 16.51302956064399
$ws.Range("G17").Value = 29.6393895984109
$ws.Range("H17").Value = 14.16003390796051
$ws.Range("I17").Value = 19.1050319021958
$ws.Range("J17").Value = 9.250815222248983
$ws.Range("N17").Value = 16.41192889064417
$ws.Range("O17").Value = 21.74404340349613
$ws.Range("B18").Value = 15.6603547649622
$ws.Range("C18").Value = 10.56762685901764
$ws.Range("D18").Value = 15.08477651647233
$ws.Range("E18").Value = 16.48984863224339
$ws.Range("G18").Value = 29.60464481530848
$ws.Range("H18").Value = 14.16686446989164
$ws.Range("I18").Value = 19.12614153121256
$ws.Range("J18").Value = 9.250666469442594
$ws.Range("N18").Value = 16.41928156459155
$ws.Range("O18").Value = 21.74654801757501
$ws.Range("B19").Value = 15.61431855086007
$ws.Range("C19").Value = 10.5286798791182
$ws.Range("D19").Value = 15.07714853185754
$ws.Range("E19").Value = 16.48205516320587
$ws.Range("G19").Value = 29.59306081693789
$ws.Range("H19").Value = 14.16922528156199
$ws.Range("I19").Value = 19.13336675919947
$ws.Range("J19").Value = 9.250640133933087
$ws.Range("N19").Value = 16.42180045914742
$ws.Range("O19").Value = 21.74749531704962
$ws.Range("B20").Value = 15.82046825589468
$ws.Range("C20").Value = 10.70281077959471
$ws.Range("D20").Value = 15.11167043761016
$ws.Range("E20").Value = 16.51734591836915
$ws.Range("G20").Value = 29.64590537908487
$ws.Range("H20").Value = 14.15879258006911
$ws.Range("I20").Value = 19.10116199392279
$ws.Range("J20").Value = 9.250854162345156
$ws.Range("N20").Value = 16.41058203946849
$ws.Range("O20").Value = 21.74362702643695
$ws.Range("B21").Value = 16.49481810235608
$ws.Range("C21").Value = 11.26783593751418
$ws.Range("D21").Value = 15.23099640261437
$ws.Range("E21").Value = 16.63967832926282
$ws.Range("G21").Value = 29.8358553452142
$ws.Range("H21").Value = 14.12709097046193
$ws.Range("I21").Value = 18.99838870669914
$ws.Range("J21").Value = 9.253241880865188
$ws.Range("N21").Value = 16.37493294193474
$ws.Range("O21").Value = 21.73754341592925
$ws.Range("B22").Value = 16.92140157676057
$ws.Range("C22").Value = 11.62200187417934
$ws.Range("D22").Value = 15.3113736767607
$ws.Range("E22").Value = 16.72233633186319
$ws.Range("G22").Value = 29.96908666214194
$ws.Range("H22").Value = 14.10887484103133
$ws.Range("I22").Value = 18.93528352317043
$ws.Range("J22").Value = 9.256042181770891
$ws.Range("N22").Value = 16.35315435482872
$ws.Range("O22").Value = 21.73872610115401
$ws.Range("B23").Value = 16.69505072564395
$ws.Range("C23").Value = 11.43436981734828
$ws.Range("D23").Value = 15.26826235200912
$ws.Range("E23").Value = 16.67797882870822
$ws.Range("G23").Value = 29.89715746071904
$ws.Range("H23").Value = 14.11836781054782
$ws.Range("I23").Value = 18.96859114785191
$ws.Range("J23").Value = 9.254434432130868
$ws.Range("N23").Value = 16.36463907421587
$ws.Range("O23").Value = 21.73762258078747
$ws.Range("B24").Value = 15.80921487289701
$ws.Range("C24").Value = 10.69332317988087
$ws.Range("D24").Value = 15.10976186492669
$ws.Range("E24").Value = 16.51539353367433
$ws.Range("G24").Value = 29.64295639346352
$ws.Range("H24").Value = 14.15935290225129
$ws.Range("I24").Value = 19.10291013658157
$ws.Range("J24").Value = 9.250836122442161
$ws.Range("N24").Value = 16.41119040769747
$ws.Range("O24").Value = 21.74381346633101
$ws.Range("B25").Value = 14.79651672213136
$ws.Range("C25").Value = 9.830399924833976
$ws.Range("D25").Value = 14.94967163685733
$ws.Range("E25").Value = 16.35224931465799
$ws.Range("G25").Value = 29.40872486113652
$ws.Range("H25").Value = 14.21427162856225
$ws.Range("I25").Value = 19.26513783562923
$ws.Range("J25").Value = 9.251091661787845
$ws.Range("N25").Value = 16.46794793573109
$ws.Range("O25").Value = 21.75995204865617
